$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("conferences")

# Insert a new row at position 2, shifting existing rows down
$ws.Rows.Item(2).Insert()

# Fill in the new conference entry
$ws.Range("A2").Value = "Presentation: The expression of elongases and desaturases shed light on the CHC plasticity of honey`nbees (\textit{Apis mellifera})"
$ws.Range("B2").Value = "September 19 2023"
$ws.Range("C2").Value = "7th meeting of the Central European Section of the IUSSE"
$ws.Range("D2").Value = "Cluj-Napoca, Romania"

$ws.Rows.Item(2).RowHeight = 60
$ws.Range("E2").Clear()

$ws.Rows.Item(5).Select()
